# Auto-generated edit script: applies 222 cell-value changes
# across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# matching the target OOXML diff for Maduin_Profits.xlsx.
#
# Each row touches columns H..N (currentAveragePrice / NQ / HQ,
# LevePriceNQ/HQ, LeveProfitNQ/HQ). Some cells are newly introduced
# (no prior value) and some are fully removed (cleared) by the edit -
# ClearContents() is used for removals so the cell disappears from
# the saved XML rather than persisting as an empty tag.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 50000
$ws.Range("J3").Value = 50000
$ws.Range("L3").Value = 50000
$ws.Range("N3").Value = -50228
$ws.Range("H11").Value = 4.6666665
$ws.Range("I11").Value = 4.6666665
$ws.Range("K11").Value = 4.6666665
$ws.Range("M11").Value = 135.3333335
$ws.Range("H33").Value = 293
$ws.Range("I33").Value = 72.59999999999999
$ws.Range("K33").Value = 72.59999999999999
$ws.Range("M33").Value = 156.4
$ws.Range("H40").Value = 1199.375
$ws.Range("J40").Value = 1198.3334
$ws.Range("L40").Value = 1198.3334
$ws.Range("N40").Value = -1548.3334
$ws.Range("H88").Value = 2270.6667
$ws.Range("I88").Value = 1432
$ws.Range("J88").Value = 2550.2222
$ws.Range("K88").Value = 1432
$ws.Range("L88").Value = 2550.2222
$ws.Range("M88").Value = -1026
$ws.Range("N88").Value = -3362.2222
$ws.Range("H91").Value = 2270.6667
$ws.Range("I91").Value = 1432
$ws.Range("J91").Value = 2550.2222
$ws.Range("K91").Value = 1432
$ws.Range("L91").Value = 2550.2222
$ws.Range("M91").Value = -28
$ws.Range("N91").Value = -5358.2222
$ws.Range("H102").Value = 50000
$ws.Range("J102").Value = 50000
$ws.Range("L102").Value = 50000
$ws.Range("N102").Value = -56490
$ws.Range("H135").Value = 451.16666
$ws.Range("I135").Value = 451.16666
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4060.49994
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1525.49994
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H102").Value = 5891
$ws.Range("I102").Value = 5891
$ws.Range("K102").Value = 5891
$ws.Range("M102").Value = -4269

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H94").Value = 3404.5
$ws.Range("I94").Value = 3404.5
$ws.Range("K94").Value = 3404.5
$ws.Range("M94").Value = -2953.5
$ws.Range("H134").Value = 390
$ws.Range("I134").Value = 390
$ws.Range("K134").Value = 1170
$ws.Range("M134").Value = 1365

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2142.2307
$ws.Range("I7").Value = 1621.2858
$ws.Range("J7").Value = 2750
$ws.Range("K7").Value = 1621.2858
$ws.Range("L7").Value = 2750
$ws.Range("M7").Value = -1508.2858
$ws.Range("N7").Value = -2976
$ws.Range("H16").Value = 2250
$ws.Range("I16").Value = 1500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1213
$ws.Range("H31").Value = 9011
$ws.Range("I31").Value = 9011
$ws.Range("K31").Value = 9011
$ws.Range("M31").Value = -8716
$ws.Range("H34").Value = 9011
$ws.Range("I34").Value = 9011
$ws.Range("K34").Value = 9011
$ws.Range("M34").Value = -8809
$ws.Range("H94").Value = 2422.5334
$ws.Range("I94").Value = 1747.6666
$ws.Range("J94").Value = 3434.8333
$ws.Range("K94").Value = 1747.6666
$ws.Range("L94").Value = 3434.8333
$ws.Range("M94").Value = -1296.6666
$ws.Range("N94").Value = -4336.8333
$ws.Range("H95").Value = 8508
$ws.Range("J95").Value = 8508
$ws.Range("L95").Value = 8508
$ws.Range("N95").Value = -14000
$ws.Range("H113").Value = 2250
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3703.2222
$ws.Range("I4").Value = 526.3333
$ws.Range("K4").Value = 1578.9999
$ws.Range("M4").Value = -1466.9999
$ws.Range("H11").Value = 2885.5715
$ws.Range("J11").Value = 3799.6
$ws.Range("L11").Value = 11398.8
$ws.Range("N11").Value = -11678.8
$ws.Range("H34").Value = 1522
$ws.Range("J34").Value = 2999.5
$ws.Range("L34").Value = 8998.5
$ws.Range("N34").Value = -9166.5
$ws.Range("H40").Value = 83.545456
$ws.Range("I40").Value = 74.8
$ws.Range("J40").Value = 90.833336
$ws.Range("K40").Value = 299.2
$ws.Range("L40").Value = 363.333344
$ws.Range("M40").Value = -230.2
$ws.Range("N40").Value = -501.333344
$ws.Range("H98").Value = 113
$ws.Range("I98").Value = 101
$ws.Range("J98").Value = 125
$ws.Range("K98").Value = 303
$ws.Range("L98").Value = 375
$ws.Range("M98").Value = 1195
$ws.Range("N98").Value = -3371
$ws.Range("H134").Value = 1034.1428
$ws.Range("I134").Value = 1034.1428
$ws.Range("K134").Value = 3102.4284
$ws.Range("M134").Value = 1967.5716
$ws.Range("H137").Value = 1165
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 2447.5
$ws.Range("I139").Value = 2447.5
$ws.Range("K139").Value = 7342.5
$ws.Range("M139").Value = -2202.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 38667.6
$ws.Range("I24").Value = 40006
$ws.Range("J24").Value = 38333
$ws.Range("K24").Value = 40006
$ws.Range("L24").Value = 38333
$ws.Range("M24").Value = -39833
$ws.Range("N24").Value = -38679
$ws.Range("H97").Value = 2368
$ws.Range("I97").Value = 2369.8333
$ws.Range("J97").Value = 2363.875
$ws.Range("K97").Value = 2369.8333
$ws.Range("L97").Value = 2363.875
$ws.Range("M97").Value = -1873.8333
$ws.Range("N97").Value = -3355.875
$ws.Range("H113").Value = 4940
$ws.Range("I113").Value = 2880.5
$ws.Range("J113").Value = 6999.5
$ws.Range("K113").Value = 2880.5
$ws.Range("L113").Value = 6999.5
$ws.Range("M113").Value = -710.5
$ws.Range("N113").Value = -11339.5
$ws.Range("H132").Value = 4222
$ws.Range("I132").Value = 3944.2
$ws.Range("K132").Value = 11832.6
$ws.Range("M132").Value = -9302.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6913.6665
$ws.Range("I16").Value = 6913.6665
$ws.Range("K16").Value = 6913.6665
$ws.Range("M16").Value = -6743.6665
$ws.Range("H46").Value = 4052.3684
$ws.Range("H55").Value = 855.55554
$ws.Range("J55").Value = 775
$ws.Range("L55").Value = 775
$ws.Range("N55").Value = -1121
$ws.Range("H82").Value = 1342.8572
$ws.Range("I82").Value = 775
$ws.Range("K82").Value = 775
$ws.Range("M82").Value = -414
$ws.Range("H85").Value = 1342.8572
$ws.Range("I85").Value = 775
$ws.Range("K85").Value = 775
$ws.Range("M85").Value = 473
$ws.Range("H93").Value = 985.4286
$ws.Range("I93").Value = 816.3333
$ws.Range("K93").Value = 816.3333
$ws.Range("M93").Value = 431.6667
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 6316.5
$ws.Range("I132").Value = 2075
$ws.Range("K132").Value = 6225
$ws.Range("M132").Value = -3695
$ws.Range("H134").Value = 65429
$ws.Range("J134").Value = 65429
$ws.Range("L134").Value = 65429
$ws.Range("N134").Value = -75569

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H54").Value = 27729.389
$ws.Range("I54").Value = 17284.75
$ws.Range("K54").Value = 17284.75
$ws.Range("M54").Value = -16764.75
$ws.Range("H105").Value = 37538
$ws.Range("J105").Value = 37538
$ws.Range("L105").Value = 37538
$ws.Range("N105").Value = -44526
$ws.Range("H122").Value = 1998.2307
$ws.Range("I122").Value = 1998.6666
$ws.Range("J122").Value = 1993
$ws.Range("K122").Value = 5995.9998
$ws.Range("L122").Value = 5979
$ws.Range("M122").Value = -3545.9998
$ws.Range("N122").Value = -10879
$ws.Range("H132").Value = 3900.5
$ws.Range("I132").Value = 3680.6
$ws.Range("K132").Value = 11041.8
$ws.Range("M132").Value = -8511.799999999999
